$d = $word.ActiveDocument

# Locate the two bullet paragraphs that need to be removed entirely:
#   "- Locker / securing data / lock "
#   "- Human Brain / storing data / recalling event"
# (searched by content rather than a hard-coded index, so the script is
# resilient to minor shifts elsewhere in the document)
$lockerPara = $null
$brainPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text

    if (($null -eq $lockerPara) -and ($text -like "*Locker*securing data*")) {
        $lockerPara = $para
    }
    elseif (($null -eq $brainPara) -and ($text -like "*Human*Brain*storing data*")) {
        $brainPara = $para
    }
}

if (($null -ne $lockerPara) -and ($null -ne $brainPara)) {
    $deleteRange = $d.Range($lockerPara.Range.Start, $brainPara.Range.End)
    $deleteRange.Delete()
}

# Change "Instructions" to "Recipe" in the remaining bullet:
#   "- Instructions / organizing data / step by step directions"
# becomes
#   "- Recipe / organizing data / step by step directions"
$d.Content.Find.Execute("Instructions", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Recipe", 2)
